$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adding updated datasheets for December 2017 - append a new weight
# measurement row (row 3) below the existing data.
$ws.Range("A3").Value = 2
$ws.Range("A3").NumberFormat = "0"

$ws.Range("B3").Value = 42736
$ws.Range("B3").NumberFormat = "m/d/yy"

$ws.Range("C3").Value = 16.7
$ws.Range("C3").NumberFormat = "0.0"

$ws.Range("D3").Value = 104.7
$ws.Range("D3").NumberFormat = "0.0"

$ws.Range("E3").Value = 230
$ws.Range("E3").NumberFormat = "0.0"

$ws.Range("F3").Formula = "=E3-E2"
$ws.Range("F3").NumberFormat = "0.0"

$ws.Range("G3").Formula = "=E3-210"
$ws.Range("G3").NumberFormat = "0.0"

# Move the active selection to reflect where the user left off editing
$ws.Range("G4").Select()
